$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. In the source sheet
# this shows up as a brand-new row 543 ("Primera" quality, recorded
# 2022-12-23) with every following record (old rows 543-606) pushed down
# by one row (new rows 544-607).
$ws.Rows("543:543").Insert()

$ws.Range("A543").Value = 3
$ws.Range("B543").Value = "Femacal de La Calera"
$ws.Range("C543").Value = "Coquimbo"
$ws.Range("D543").Value = 44918
$ws.Range("E543").Value = 5
$ws.Range("F543").Value = 100112037
$ws.Range("G543").Value = "Cebollín"
$ws.Range("H543").Value = "Sin especificar"
$ws.Range("I543").Value = "Primera"
$ws.Range("J543").Value = 270
$ws.Range("K543").Value = 3500
$ws.Range("L543").Value = 4000
$ws.Range("M543").Value = 3778
$ws.Range("N543").Value = "$/paquete 36 unidades"
$ws.Range("O543").Value = "Provincia de Quillota"
$ws.Range("P543").Value = 105
$ws.Range("Q543").Value = 36
$ws.Range("R543").Value = "Hortaliza"
